$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.956.53"
$ws.Range('D3').Value = "'1.638.15"
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'212.32"
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = "'23.37"
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  -2.67%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = "'1.870.61"
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = "'1.637.51"
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = "'0.568"
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').Value = "'65.41"
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = "'27.962.63"
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = "'232.38"
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('D23').Value = "'4.36"
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = "'2.07"
$ws.Range('E24').Value = '  -3.58%  '
$ws.Range('D25').Value = "'153.18"
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('D26').Value = "'6.97"
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.111"
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'15.63"
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').Value = "'3.39"
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').Value = "'1.403.56"
$ws.Range('E33').Value = '  -3.64%  '
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('D36').Value = "'2.36"
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').Value = "'0.926"
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('D40').Value = "'0.875"
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = "'67.04"
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('E44').Value = '  +2.88%  '
$ws.Range('D45').Value = "'1.81"
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').Value = "'1.779.68"
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').Value = "'87.98"
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = "'7.56"
$ws.Range('E51').Value = '  -2.14%  '
